$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------------
# 1. Update the "Status" text from "Handed back: in sync with en-US" to
#    "Ready for handoff" for the row that survives (the 653de435... row).
# ---------------------------------------------------------------------------
$wsOverview.Range("B2").Value2 = "Ready for handoff"
$wsOverview.Range("C2").Value2 = "Ready for handoff"
$wsZhCn.Range("B2").Value2 = "Ready for handoff"
$wsDeDe.Range("B2").Value2 = "Ready for handoff"

# ---------------------------------------------------------------------------
# 2. Update the "Latest Handoff Datetime" values for the 653de435... entries.
# ---------------------------------------------------------------------------
$wsZhCn.Range("D2").Value2 = "2016-02-18 09:54:23"
$wsDeDe.Range("D2").Value2 = "2016-02-18 09:54:33"

# ---------------------------------------------------------------------------
# 3. Remove the whole "a634b5f3-a252-4698-b996-c9ad1c439b66..." row (row 3)
#    from every sheet. This shifts the ".localization-config" row up from
#    row 4 to row 3, exactly as described in the diff.
# ---------------------------------------------------------------------------
$wsOverview.Rows.Item(3).Delete()
$wsZhCn.Rows.Item(3).Delete()
$wsDeDe.Rows.Item(3).Delete()

# ---------------------------------------------------------------------------
# 4. Rebuild the hyperlinks on every sheet: the row delete above does not
#    keep the hyperlink collection in sync, so clear it out and re-create
#    only the links that should remain, pointing at the correct (shifted)
#    cells with their original target URLs / display text.
# ---------------------------------------------------------------------------
$wsOverview.Range("A1:C3").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/09eab70e082b39555f06f9615a75860a53d3acdc/e2e/653de435-9040-4ca6-864d-6e5c29891627.md", "", "", "653de435-9040-4ca6-864d-6e5c29891627.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/09eab70e082b39555f06f9615a75860a53d3acdc/.localization-config", "", "", ".localization-config")

$wsZhCn.Range("A1:I3").Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/09eab70e082b39555f06f9615a75860a53d3acdc/e2e/653de435-9040-4ca6-864d-6e5c29891627.md", "", "", "653de435-9040-4ca6-864d-6e5c29891627.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2067413323a0893eee1f0756ac9235d38d1752b8/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/653de435-9040-4ca6-864d-6e5c29891627.dcb53bc46cef1baa3842aff3ceef28a8e9bffd2f.zh-cn.xlf", "", "", "653de435-9040-4ca6-864d-6e5c29891627.dcb53bc46cef1baa3842aff3ceef28a8e9bffd2f.zh-cn.xlf")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/32178978a4ee0bd9bae4df385466ef1fa24f23cf/e2e/653de435-9040-4ca6-864d-6e5c29891627.md", "", "", "653de435-9040-4ca6-864d-6e5c29891627.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/70fa8d87806aba64f9c5088bfd5bda46297eb980/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/653de435-9040-4ca6-864d-6e5c29891627.dcb53bc46cef1baa3842aff3ceef28a8e9bffd2f.zh-cn.xlf", "", "", "653de435-9040-4ca6-864d-6e5c29891627.dcb53bc46cef1baa3842aff3ceef28a8e9bffd2f.zh-cn.xlf")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/09eab70e082b39555f06f9615a75860a53d3acdc/.localization-config", "", "", ".localization-config")

$wsDeDe.Range("A1:I3").Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/09eab70e082b39555f06f9615a75860a53d3acdc/e2e/653de435-9040-4ca6-864d-6e5c29891627.md", "", "", "653de435-9040-4ca6-864d-6e5c29891627.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/054b5e0f95c4da417ff3b2a9e824a5df884385ac/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/653de435-9040-4ca6-864d-6e5c29891627.dcb53bc46cef1baa3842aff3ceef28a8e9bffd2f.de-de.xlf", "", "", "653de435-9040-4ca6-864d-6e5c29891627.dcb53bc46cef1baa3842aff3ceef28a8e9bffd2f.de-de.xlf")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/7fe05ba958f05c6cdbb6ba4172ff5c91032917c5/e2e/653de435-9040-4ca6-864d-6e5c29891627.md", "", "", "653de435-9040-4ca6-864d-6e5c29891627.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/c635ff0284f6fdcd49d5545c817c3230f3209950/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/653de435-9040-4ca6-864d-6e5c29891627.dcb53bc46cef1baa3842aff3ceef28a8e9bffd2f.de-de.xlf", "", "", "653de435-9040-4ca6-864d-6e5c29891627.dcb53bc46cef1baa3842aff3ceef28a8e9bffd2f.de-de.xlf")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/09eab70e082b39555f06f9615a75860a53d3acdc/.localization-config", "", "", ".localization-config")
